$wb = $excel.ActiveWorkbook

# --- Update existing "Lead Tracking" sheet: append two more rows for Piotr P ---
$ws1 = $wb.Worksheets.Item("Lead Tracking")

$ws1.Range("A4").Value = "Piotr P"
$ws1.Range("B4").Value = "Technical Leader Allegro"
$ws1.Range("C4").Value = "500+"

$ws1.Range("A5").Value = "Piotr P"
$ws1.Range("B5").Value = "Technical Leader Allegro"
$ws1.Range("C5").Value = "500+"

# --- Add a new sheet "lead tracking1" right after "Lead Tracking" ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "lead tracking1"

$newSheet.Range("A1").Value = "name"
$newSheet.Range("B1").Value = "designation"
$newSheet.Range("C1").Value = "total connection"

$newSheet.Range("A2").Value = "Piotr P"
$newSheet.Range("B2").Value = "Technical Leader Allegro"
$newSheet.Range("C2").Value = "500+"
